$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "BVT": add a new BVT checklist row (row 7 - "Negative Support")
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("BVT")

$ws1.Range("A7").Value = 6
$ws1.Range("B7").Value = "Negative Support"
$ws1.Range("C7").Value = "Drag columns with negative values"
$ws1.Range("D7").Value = "1. Drag 'Organization' column in 'Category " + [char]10 + "Data' field" + [char]10 + "2.Drag 'Quarter 1 Growth/Fall' in 'Measure Data'" + [char]10 + "3. Drag 'Quarter 2 Growth/Fall' in 'Measure Data'"
$ws1.Range("E7").Value = "Plot should render for negative values also and data labels should also appear for the respective bars."

# New row mirrors the wrap-text styling used by the rest of the D/E columns
$ws1.Range("D7:E7").WrapText = $true

# Row heights were recalculated (file re-saved from a newer Excel build)
$ws1.Rows.Item(2).RowHeight = 75
$ws1.Rows.Item(3).RowHeight = 45
$ws1.Rows.Item(4).RowHeight = 45
$ws1.Rows.Item(5).RowHeight = 90
$ws1.Rows.Item(6).RowHeight = 60
$ws1.Rows.Item(7).RowHeight = 90

# Selection moved to the newly added cell
$null = $ws1.Range("E7").Select()

# ----------------------------------------------------------------------
# Sheet "Checklist": row heights recalculated the same way
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Checklist")

$ws2.Rows.Item(1).RowHeight = 15.75
$ws2.Rows.Item(2).RowHeight = 30.75
$ws2.Rows.Item(3).RowHeight = 15.75
$ws2.Rows.Item(4).RowHeight = 15.75
$ws2.Rows.Item(5).RowHeight = 15.75
$ws2.Rows.Item(6).RowHeight = 45.75
$ws2.Rows.Item(7).RowHeight = 30.75
$ws2.Rows.Item(8).RowHeight = 30.75
$ws2.Rows.Item(9).RowHeight = 15.75
$ws2.Rows.Item(10).RowHeight = 15.75
$ws2.Rows.Item(11).RowHeight = 15.75
$ws2.Rows.Item(12).RowHeight = 15.75
$ws2.Rows.Item(13).RowHeight = 15.75
$ws2.Rows.Item(14).RowHeight = 30.75
$ws2.Rows.Item(15).RowHeight = 15.75
$ws2.Rows.Item(16).RowHeight = 30.75
$ws2.Rows.Item(19).RowHeight = 30.75
$ws2.Rows.Item(20).RowHeight = 30.75
$ws2.Rows.Item(25).RowHeight = 15.75
$ws2.Rows.Item(26).RowHeight = 15.75
$ws2.Rows.Item(27).RowHeight = 15.75

# Restore BVT as the active sheet (matches tabSelected="1" in the source)
$ws1.Activate()
